# Auto-generated Excel COM-interop script to apply cryptos list update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking text values in Price/Volume columns are kept as Text
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "79.759.84"
$ws.Range("E2").Value = "  +4.87%  "
$ws.Range("D3").Value = "3.213.35"
$ws.Range("E3").Value = "  +6.18%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "211.24"
$ws.Range("E5").Value = "  +7.41%  "
$ws.Range("D6").Value = "638.22"
$ws.Range("E6").Value = "  +3.11%  "
$ws.Range("D7").Value = "0.263"
$ws.Range("E7").Value = "  +28.13%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "0.604"
$ws.Range("E9").Value = "  +10.01%  "
$ws.Range("D10").Value = "3.211.39"
$ws.Range("E10").Value = "  +6.17%  "
$ws.Range("D11").Value = "0.608"
$ws.Range("E11").Value = "  +38.27%  "
$ws.Range("D12").Value = "0.0000268"
$ws.Range("E12").Value = "  +40.08%  "
$ws.Range("E13").Value = "  +3.54%  "
$ws.Range("D14").Value = "5.44"
$ws.Range("E14").Value = "  +4.06%  "
$ws.Range("D15").Value = "3.802.06"
$ws.Range("E15").Value = "  +6.03%  "
$ws.Range("D16").Value = "32.76"
$ws.Range("E16").Value = "  +13.63%  "
$ws.Range("D17").Value = "79.554.59"
$ws.Range("E17").Value = "  +4.65%  "
$ws.Range("D18").Value = "3.203.75"
$ws.Range("E18").Value = "  +5.92%  "
$ws.Range("D19").Value = "14.67"
$ws.Range("E19").Value = "  +9.28%  "
$ws.Range("D20").Value = "9.44"
$ws.Range("E20").Value = "  +5.96%  "
$ws.Range("B21").Value = "SuiNetwork"
$ws.Range("C21").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D21").Value = "3.02"
$ws.Range("E21").Value = "  +28.13%  "
$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").Value = "448.99"
$ws.Range("E22").Value = "  +18.32%  "
$ws.Range("D23").Value = "5.30"
$ws.Range("E23").Value = "  +21.14%  "
$ws.Range("D24").Value = "4.86"
$ws.Range("E24").Value = "  +12.87%  "
$ws.Range("D25").Value = "3.366.32"
$ws.Range("E25").Value = "  +5.77%  "
$ws.Range("D26").Value = "77.88"
$ws.Range("E26").Value = "  +7.54%  "
$ws.Range("D27").Value = "10.95"
$ws.Range("E27").Value = "  +12.23%  "
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  +0.10%  "
$ws.Range("D29").Value = "0.0000126"
$ws.Range("E29").Value = "  +17.71%  "
$ws.Range("D30").Value = "9.25"
$ws.Range("E30").Value = "  +12.47%  "
$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("D32").Value = "566.97"
$ws.Range("E32").Value = "  +15.48%  "
$ws.Range("D33").Value = "1.53"
$ws.Range("E33").Value = "  +10.11%  "
$ws.Range("D34").Value = "0.157"
$ws.Range("E34").Value = "  +28.49%  "
$ws.Range("E35").Value = "  +6.25%  "
$ws.Range("D36").Value = "23.24"
$ws.Range("E36").Value = "  +13.26%  "
$ws.Range("D37").Value = "0.123"
$ws.Range("E37").Value = "  +20.41%  "
$ws.Range("E38").Value = "  -0.05%  "
$ws.Range("D39").Value = "0.416"
$ws.Range("E39").Value = "  +10.20%  "
$ws.Range("D40").Value = "163.42"
$ws.Range("E40").Value = "  +0.86%  "
$ws.Range("D41").Value = "20.29"
$ws.Range("E41").Value = "  +1.38%  "
$ws.Range("E42").Value = "  +12.66%  "
$ws.Range("D43").Value = "192.67"
$ws.Range("E43").Value = "  +1.33%  "
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("E45").Value = "  +11.98%  "
$ws.Range("B46").Value = "dogwifhat"
$ws.Range("C46").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D46").Value = "2.73"
$ws.Range("E46").Value = "  +13.42%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").Value = "0.802"
$ws.Range("E47").Value = "  +4.66%  "
$ws.Range("D48").Value = "1.35"
$ws.Range("E48").Value = "  +8.61%  "
$ws.Range("D49").Value = "43.23"
$ws.Range("E49").Value = "  +4.62%  "
$ws.Range("B50").Value = "Filecoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D50").Value = "4.33"
$ws.Range("E50").Value = "  +12.43%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").Value = "25.98"
$ws.Range("E51").Value = "  +16.81%  "
